$wb = $excel.ActiveWorkbook

# Update "想去人数" (F6, F9) counts by +1 on both the "展览" and "全部类型" sheets
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F6").Value = 436
    $ws.Range("F9").Value = 557
}
